# Update countries & provincias Spain
# Refresh the COVID country-stats sheet: bump the "last updated" timestamp,
# push new per-country totals, and re-sort a handful of countries whose
# case counts now place them in a different rank position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Septiembre de 2020 a las 05:51"

# --- Kazajistan (row 35) - stat refresh only ---------------------------
$ws.Range("B35").Value = 107262
$ws.Range("C35").Value = 63
$ws.Range("D35").Value = 101877
$ws.Range("E35").Value = 3714

# --- Belgica now outranks Marruecos: swap rows 38/39 --------------------
$ws.Range("A38").Value = "Belgica"
$ws.Range("B38").Value = 100748
$ws.Range("C38").Value = 1099
$ws.Range("D38").Value = 18945
$ws.Range("E38").Value = 71859
$ws.Range("G38").Value = 7
$ws.Range("H38").Value = 9944

$ws.Range("A39").Value = "Marruecos"
$ws.Range("B39").Value = 99816
$ws.Range("D39").Value = 79008
$ws.Range("E39").Value = 19013
$ws.Range("H39").Value = 1795

# --- Honduras (row 50) - stat refresh only ------------------------------
$ws.Range("B50").Value = 71143
$ws.Range("C50").Value = 532
$ws.Range("D50").Value = 21810
$ws.Range("E50").Value = 47167
$ws.Range("G50").Value = 20
$ws.Range("H50").Value = 2166

# --- Birmania re-enters the ranking at row 115, pushing rows 115-121 ----
# --- (Suazilandia..Nicaragua) down one slot; old row 122 data retired --
$ws.Range("A115").Value = "Birmania"
$ws.Range("B115").Value = 5263
$ws.Range("C115").Value = 393
$ws.Range("D115").Value = 1188
$ws.Range("E115").Value = 3994
$ws.Range("H115").Value = 81

$ws.Range("A116").Value = "Suazilandia"
$ws.Range("B116").Value = 5245
$ws.Range("D116").Value = 4571
$ws.Range("E116").Value = 570
$ws.Range("H116").Value = 104

$ws.Range("A117").Value = "Cabo Verde"
$ws.Range("B117").Value = 5186
$ws.Range("D117").Value = 4581
$ws.Range("E117").Value = 555
$ws.Range("H117").Value = 50

$ws.Range("A118").Value = "Cuba"
$ws.Range("B118").Value = 5055
$ws.Range("D118").Value = 4284
$ws.Range("E118").Value = 658
$ws.Range("H118").Value = 113

$ws.Range("A119").Value = "Hong Kong"
$ws.Range("B119").Value = 5010
$ws.Range("D119").Value = 4707
$ws.Range("E119").Value = 200
$ws.Range("H119").Value = 103

$ws.Range("A120").Value = "Guinea Ecuatorial"
$ws.Range("B120").Value = 5002
$ws.Range("D120").Value = 4509
$ws.Range("E120").Value = 410
$ws.Range("H120").Value = 83

$ws.Range("A121").Value = "Congo"
$ws.Range("B121").Value = 4986
$ws.Range("D121").Value = 3887
$ws.Range("E121").Value = 1010
$ws.Range("H121").Value = 89

$ws.Range("A122").Value = "Nicaragua"
$ws.Range("B122").Value = 4961
$ws.Range("D122").Value = 2913
$ws.Range("E122").Value = 1901
$ws.Range("H122").Value = 147

# --- Belice (row 158) - stat refresh only -------------------------------
$ws.Range("B158").Value = 1606
$ws.Range("C158").Value = 16
$ws.Range("D158").Value = 876
$ws.Range("E158").Value = 710

# --- Polinesia Francesa (row 164) - stat refresh only -------------------
$ws.Range("B164").Value = 1271
$ws.Range("D164").Value = 1028

# --- Islas Turcas y Caicos (row 172) - stat refresh only ----------------
$ws.Range("B172").Value = 668
$ws.Range("C172").Value = 1
$ws.Range("D172").Value = 572
$ws.Range("E172").Value = 91

# --- San Martin (Parte Holandesa) (row 173) - stat refresh only ---------
$ws.Range("B173").Value = 584
$ws.Range("C173").Value = 10
$ws.Range("E173").Value = 76

# --- Mongolia (row 185) - stat refresh only ------------------------------
$ws.Range("B185").Value = 312
$ws.Range("C185").Value = 1
$ws.Range("E185").Value = 10

# --- Butan (row 187) - stat refresh only ---------------------------------
$ws.Range("B187").Value = 259
$ws.Range("C187").Value = 1
$ws.Range("D187").Value = 190
$ws.Range("E187").Value = 69

# --- Timor Oriental now ties/outranks Santa Lucia: swap row labels ------
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"

# --- Islas Malvinas now outranks Montserrat: swap rows 214/215 ----------
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
